# Applies the "fixed spelling errors on poster" commit:
#  - bumps the cached Date placeholder text (4/7/2022 -> 4/12/2022) on the
#    slide master, every slide layout, and the notes master
#  - fixes a handful of small wording/punctuation issues on slide 1

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $find inside $shape's text with
# $replace, touching only that sub-range so the rest of the run(s) keep their
# original formatting untouched.
# ---------------------------------------------------------------------------
function Replace-InShapeText($shape, $find, $replace) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $find.Length)
        $sub.Text = $replace
    }
}

function Get-ShapeByName($container, $name) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

function Set-DateFigureText($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# ---------------------------------------------------------------------------
# 1. Date placeholder: 4/7/2022 -> 4/12/2022
#    (slide master + every layout; notes master goes through HeadersFooters
#    since touching its Shapes collection directly doesn't stick)
# ---------------------------------------------------------------------------
$newDate = "4/12/2022"

$master = $p.SlideMaster
Set-DateFigureText $master $newDate

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    Set-DateFigureText ($master.CustomLayouts.Item($j)) $newDate
}

$notesDt = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDt.Text = $newDate

# ---------------------------------------------------------------------------
# 2. Slide 1 text fixes
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$tb15 = Get-ShapeByName $slide "TextBox 15"
Replace-InShapeText $tb15 "web based application" "web-based application"

$tb51 = Get-ShapeByName $slide "TextBox 51"
Replace-InShapeText $tb51 "reaching the database, or searching" "reaching the database or searching"
Replace-InShapeText $tb51 "(ERD shown in Figure 2) will include four tables" "(ERD shown in Figure 2) includes four tables"
Replace-InShapeText $tb51 "The User table will have 3 classes" "The User table has 3 classes"

$tb53 = Get-ShapeByName $slide "TextBox 53"
Replace-InShapeText $tb53 "would like to add, but will not" "would like to add but will not"

$tb61 = Get-ShapeByName $slide "TextBox 61"
Replace-InShapeText $tb61 "their own shades there will be" "their own shades, there will be"

# "When users save a palette it will automatically" -> split into three runs
# with a comma inserted: "...save " / "a palette, " / "it will automatically..."
$tr61 = $tb61.TextFrame.TextRange
$full61 = $tr61.Text
$marker = "When users save "
$startIdx = $full61.IndexOf($marker) + $marker.Length
$afterIdx = $full61.IndexOf("it will automatically", $startIdx)
$segLen = $afterIdx - $startIdx
$seg = $tr61.Characters($startIdx + 1, $segLen)
$seg.Text = "a palette, "
